$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73, shifting existing rows 73-86 down to 74-87
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record
$ws.Cells.Item(73, 1).Value = 5
$ws.Cells.Item(73, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(73, 3).Value = "Maule"
$ws.Cells.Item(73, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(73, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(73, 5).Value = 7
$ws.Cells.Item(73, 6).Value = 100112022
$ws.Cells.Item(73, 7).Value = "Arveja Verde"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 300
$ws.Cells.Item(73, 11).Value = 15000
$ws.Cells.Item(73, 12).Value = 15000
$ws.Cells.Item(73, 13).Value = 15000
$ws.Cells.Item(73, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(73, 15).Value = "Carahue"
$ws.Cells.Item(73, 16).Value = 600
$ws.Cells.Item(73, 17).Value = 25
$ws.Cells.Item(73, 18).Value = "Hortaliza"
